# Update a set of numeric results in Sheet1 (result_data_RandomForest.xlsx)
# These values correspond to recomputed algorithm results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 16.3891
$ws.Range("C7").Value = -13.28699999999999
$ws.Range("B10").Value = 5.581299999999998
$ws.Range("B12").Value = 4.992799999999999
$ws.Range("C15").Value = -14.18709999999999
$ws.Range("B18").Value = 6.648999999999993
$ws.Range("E18").Value = 18.0276
$ws.Range("E19").Value = 16.5417
$ws.Range("C20").Value = -12.0556
$ws.Range("E27").Value = 16.54599999999999
$ws.Range("C29").Value = -11.4575
$ws.Range("C30").Value = -12.76129999999999
$ws.Range("C31").Value = -12.78379999999999
$ws.Range("B37").Value = 8.964
$ws.Range("C40").Value = -12.6933
$ws.Range("E42").Value = 16.563
$ws.Range("E44").Value = 16.47839999999999
$ws.Range("E47").Value = 16.52869999999999
$ws.Range("B55").Value = 6.266599999999994
$ws.Range("E58").Value = 16.445
$ws.Range("B68").Value = 5.899900000000001
$ws.Range("C68").Value = -12.1249
$ws.Range("E73").Value = 17.35050000000001
$ws.Range("C76").Value = -12.17780000000001
$ws.Range("B77").Value = 8.773300000000008
$ws.Range("B78").Value = 9.3781
$ws.Range("C87").Value = -13.99499999999999
$ws.Range("C88").Value = -13.07339999999999
$ws.Range("E95").Value = 18.18690000000002
$ws.Range("C96").Value = -12.604
$ws.Range("C98").Value = -12.0616
$ws.Range("C101").Value = -12.2627
$ws.Range("E101").Value = 16.66220000000002
$ws.Range("C102").Value = -12.6751
